$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the adata_path column (E) to point at the new "data_balanced" folder
# instead of the old "data" folder.
$ws.Cells.Replace("./data/", "./data_balanced/", 2)

# Update the dose_subset column (F): rows that previously used a 1000-cell
# downsampled subset now use 10000.
$fRange = $ws.Range("F2:F141")
$fRange.Replace(1000, 10000, 1)
